$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Blackbox sheet: insert a new "Actual Results" column just before the
# existing "Pass/Fail" column, and fill in its values.
# ---------------------------------------------------------------------
$bb = $wb.Worksheets.Item("Blackbox")
$bb.Activate()
$bb.Range("F1").EntireColumn.Insert()

$bb.Cells.Item(2, 6).Value2 = "Actual`nResults"
$bb.Cells.Item(3, 6).Value2 = "MATCHES EXPECTED RESULTS"
$bb.Cells.Item(4, 6).Value2 = "MATCHES EXPECTED RESULTS"
$bb.Cells.Item(5, 6).Value2 = "MATCHES EXPECTED RESULTS"
$bb.Cells.Item(6, 6).Value2 = "MATCHES EXPECTED RESULTS"
$bb.Cells.Item(7, 6).Value2 = "MATCHES EXPECTED RESULTS"
$bb.Cells.Item(8, 6).Value2 = "MATCHES EXPECTED RESULTS"
$bb.Cells.Item(9, 6).Value2 = "MATCHES EXPECTED RESULTS"

$bb.Range("F8").Select()

# ---------------------------------------------------------------------
# Whitebox sheet: same new column, header differs slightly (no newline).
# ---------------------------------------------------------------------
$wbx = $wb.Worksheets.Item("Whitebox")
$wbx.Activate()
$wbx.Range("F1").EntireColumn.Insert()

$wbx.Cells.Item(2, 6).Value2 = "Actual Results"
$wbx.Cells.Item(3, 6).Value2 = "MATCHES EXPECTED RESULTS"
$wbx.Cells.Item(4, 6).Value2 = "MATCHES EXPECTED RESULTS"
$wbx.Cells.Item(5, 6).Value2 = "MATCHES EXPECTED RESULTS"
$wbx.Cells.Item(6, 6).Value2 = "MATCHES EXPECTED RESULTS"
$wbx.Cells.Item(7, 6).Value2 = "MATCHES EXPECTED RESULTS"
$wbx.Cells.Item(8, 6).Value2 = "MATCHES EXPECTED RESULTS"
$wbx.Cells.Item(9, 6).Value2 = "MATCHES EXPECTED RESULTS"
$wbx.Cells.Item(10, 6).Value2 = "MATCHES EXPECTED RESULTS"

$wbx.Range("F3").Select()

# ---------------------------------------------------------------------
# Description sheet becomes the active / selected tab again.
# ---------------------------------------------------------------------
$desc = $wb.Worksheets.Item("Description")
$desc.Activate()
$desc.Range("N9").Select()
